$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D (shifts old D:K data to F:M), limited to the used data range
$ws.Range("D5:E102").Insert(-4161)

# Copy number formats/styles from the (now-shifted) F:G columns into the new D:E columns
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)

# Populate the two new quarters (2018-12-31 and 2018-09-30) of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 145500
$ws.Range("E8").Value = 145700
$ws.Range("D9").Value = 25500
$ws.Range("E9").Value = 27200
$ws.Range("D10").Value = 120000
$ws.Range("E10").Value = 118500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 2300
$ws.Range("E15").Value = 2400
$ws.Range("D17").Value = 125300
$ws.Range("E17").Value = 126400
$ws.Range("D18").Value = 20200
$ws.Range("E18").Value = 19300
$ws.Range("D20").Value = 1200
$ws.Range("E20").Value = 1000
$ws.Range("D21").Value = 23700
$ws.Range("E21").Value = 22600
$ws.Range("D22").Value = 400
$ws.Range("E22").Value = 100
$ws.Range("D23").Value = 21000
$ws.Range("E23").Value = 20200
$ws.Range("D24").Value = 7000
$ws.Range("E24").Value = 5100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 14000
$ws.Range("E26").Value = 15100
$ws.Range("D27").Value = 14000
$ws.Range("E27").Value = 15100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 100
$ws.Range("E29").Value = -200
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1200
$ws.Range("E32").Value = -1000
$ws.Range("D33").Value = 14100
$ws.Range("E33").Value = 14900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 14100
$ws.Range("E35").Value = 14900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 32400
$ws.Range("E41").Value = 32000
$ws.Range("D42").Value = 196400
$ws.Range("E42").Value = 156300
$ws.Range("D43").Value = 31300
$ws.Range("E43").Value = 36500
$ws.Range("D44").Value = 800
$ws.Range("E44").Value = 900
$ws.Range("D45").Value = 8500
$ws.Range("E45").Value = 14400
$ws.Range("D46").Value = 269400
$ws.Range("E46").Value = 240100
$ws.Range("D47").Value = 900
$ws.Range("E47").Value = 2200
$ws.Range("D48").Value = 30000
$ws.Range("E48").Value = 30000
$ws.Range("D49").Value = 95300
$ws.Range("E49").Value = 95300
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 86800
$ws.Range("E52").Value = 93800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 482500
$ws.Range("E54").Value = 461400
$ws.Range("D57").Value = 9200
$ws.Range("E57").Value = 12600
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 87900
$ws.Range("E59").Value = 79000
$ws.Range("D60").Value = 97100
$ws.Range("E60").Value = 91600
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 30200
$ws.Range("E62").Value = 30100
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 127300
$ws.Range("E66").Value = 121700
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -52900
$ws.Range("E72").Value = -67000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 355200
$ws.Range("E76").Value = 339600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 14100
$ws.Range("E81").Value = 14900
$ws.Range("D83").Value = 2300
$ws.Range("E83").Value = 2400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 38700
$ws.Range("E89").Value = 3500
$ws.Range("D91").Value = -2800
$ws.Range("E91").Value = -1200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -38400
$ws.Range("E94").Value = -5500
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 300
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 400
$ws.Range("E102").Value = -1600

Write-Host "Done"
